# Updated cryptos list (price + 1h volume change) per the GitHub Actions refresh.
# Cells in column D that look like plain numbers are written with a leading
# apostrophe so Excel stores them as text (matching the original inlineStr
# cells) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.297.64"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.928.59"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'0.7541"
$ws.Range("E5").Value = "  +4.46%  "

$ws.Range("D6").Value = "'242.99"
$ws.Range("E6").Value = "  -2.73%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'28.23"
$ws.Range("E8").Value = "  +2.48%  "

$ws.Range("D9").Value = "'0.3167"
$ws.Range("E9").Value = "  -1.47%  "

$ws.Range("D10").Value = "'0.07017"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.7768"
$ws.Range("E11").Value = "  -1.89%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.08012"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "1.933.38"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "'5.365"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").Value = "'93.33"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("D16").Value = "'14.43"
$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("D17").Value = "30.312.40"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").Value = "'252.98"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "'0.000007920"
$ws.Range("E19").Value = "  -1.55%  "

$ws.Range("D20").Value = "'5.797"
$ws.Range("E20").Value = "  +1.11%  "

$ws.Range("D21").Value = "2.173.83"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'6.692"
$ws.Range("E24").Value = "  -2.04%  "

$ws.Range("D25").Value = "'9.470"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("D26").Value = "'164.70"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "'19.05"
$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("D28").Value = "'0.1332"
$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("D29").Value = "'2.197"
$ws.Range("E29").Value = "  -3.82%  "

$ws.Range("D30").Value = "'1.363"
$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("D31").Value = "'1.509"
$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").Value = "'4.396"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").Value = "'4.116"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").Value = "'1.330"
$ws.Range("E34").Value = "  +5.70%  "

$ws.Range("D35").Value = "'0.05162"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").Value = "'0.7530"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "'0.01953"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'2.793"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.509"
$ws.Range("E40").Value = "  +2.02%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'77.70"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("D42").Value = "'0.4481"
$ws.Range("E42").Value = "  -0.48%  "

$ws.Range("D43").Value = "'1.967"
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'0.8340"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").Value = "'101.40"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.866"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.558"
$ws.Range("E48").Value = "  +1.52%  "

$ws.Range("D49").Value = "'987.00"
$ws.Range("E49").Value = "  +6.67%  "

$ws.Range("D50").Value = "'37.52"
$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1193"
$ws.Range("E51").Value = "  +5.46%  "
